$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 1-3 are blank placeholder rows sitting above the table header (which
# currently lives on row 4, with data in rows 5-146). Deleting them shifts
# the whole table up by three rows, so the header ends up on row 1 and the
# data occupies rows 2-143, matching the new dimension A1:E143.
$ws.Rows("1:3").Delete()

# Move the active selection/cursor to G9 (previously D119, before the
# scroll position reset that comes from no longer needing to view far down
# the now-shorter sheet).
$ws.Range("G9").Select()
